$d = $word.ActiveDocument

# The date line currently reads "Karawaci, Juni 2022" (the "J" of "Juni"
# starts its own run right after the space that follows the comma). The
# edit adds the missing day-of-month, turning it into
# "Karawaci, 20 Juni 2022" by inserting "20 " right before "Juni".
$rng = $d.Content
$found = $rng.Find.Execute("Karawaci, ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of "Karawaci, " - i.e. the insertion point that
    # sits immediately before the "J" of "Juni" - and insert the day there.
    $rng.Collapse(0)
    $rng.InsertBefore("20 ")
}
